$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5156313333333333
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 1.498861050651777
$ws.Range("R2").Value = 13.489749455866
$ws.Range("S2").Value = 0.005520525738044089
$ws.Range("T2").Value = 0.005624540846623205

# Row 3
$ws.Range("G3").Value = 0.5156313333333333
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 95.8057167975891
$ws.Range("R3").Value = 862.251451178302
$ws.Range("S3").Value = 0.3528665483720876
$ws.Range("T3").Value = 0.3595150912979765

# Row 4
$ws.Range("G4").Value = 0.5156313333333333
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 70.67846790207155
$ws.Range("R4").Value = 636.106211118644
$ws.Range("S4").Value = 0.2603191943704447
$ws.Range("T4").Value = 0.2652240042658267

# Row 5
$ws.Range("G5").Value = 0.5156313333333333
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 15.06300227181033
$ws.Range("R5").Value = 90.378013630862
$ws.Range("S5").Value = 0.05547925319534149
$ws.Range("T5").Value = 0.03768304451958546

# Row 6
$ws.Range("G6").Value = 0.5156313333333333
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 88.46089203660088
$ws.Range("R6").Value = 796.1480283294079
$ws.Range("S6").Value = 0.3258144783240821
$ws.Range("T6").Value = 0.331953319069988
